$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1копейка ")
Write-Host $ws.Name
